$d = $word.ActiveDocument

# Update the date heading in the first paragraph.
$d.Content.Find.Execute("2023-12-26 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-12-27 Wednesday", 2)

# Update the multiplication problems in the table. Each content row holds five
# cells; several new values collide with *other* old values elsewhere in the
# table, so we must set each cell's text positionally (by row/column) rather
# than via a global Find & Replace, which could clobber cells out of order.
$t = $d.Tables.Item(1)

$values = @(
    @("54×69=", "51×78=", "77×75=", "45×12=", "70×34="),
    @("82×77=", "17×77=", "89×73=", "54×31=", "94×91="),
    @("99×74=", "94×39=", "73×65=", "34×60=", "17×67="),
    @("28×24=", "21×56=", "66×60=", "11×67=", "99×75="),
    @("67×45=", "54×53=", "52×98=", "90×30=", "84×94=")
)

$rowIndexes = @(1, 5, 10, 15, 20)

for ($i = 0; $i -lt $rowIndexes.Length; $i++) {
    $row = $t.Rows.Item($rowIndexes[$i])
    $rowValues = $values[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $row.Cells.Item($col).Range.Text = $rowValues[$col - 1]
    }
}
